$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (A6) used to carry the "last row" date-only style; now it should
# match the regular rows (same style as A2:A5).
$ws.Range("A6").NumberFormat = $ws.Range("A2").NumberFormat

# New row 7: date + profit values, with A7 taking over the "last row" style
# that A6 previously had.
$ws.Range("A7").Value = 44516
$ws.Range("A7").NumberFormat = "YYYY-MM-DD"
$ws.Range("B7").Value = -1647.2
